# Update price-list date and two prices on "Hoja1"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds the list date (serial date value), move it one month forward
$ws.Range("A1").Value = 45436

# D35 / D36 hold prices, bump them up
$ws.Range("D35").Value = 32935
$ws.Range("D36").Value = 7128
